# "upgrade left table until javakheti" — add the 2023 data column (K) to the
# Akhaltsikhe average-monthly-remuneration table, matching the formatting
# already used by the other year columns (B:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring column K's formatting in line with the existing year column (J) --
# same number format / alignment / borders used for the header row and the
# three data rows -- before writing the new values into it.
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New year header
$ws.Range("K3").Value = 2023

# New figures for the three data rows (Average / Women / Men)
$ws.Range("K4").Value = 1025.5
$ws.Range("K5").Value = 543.8
$ws.Range("K6").Value = 1299.2
